$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "28.666.49"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -1.96%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.802.91"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  -1.52%  "

$ws.Range("E4").Value = "  +0.20%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "231.35"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -2.35%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.5960"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -1.95%  "

$ws.Range("E7").Value = "  +0.21%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.2772"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -1.63%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.06841"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -3.91%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "23.34"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -2.64%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.07529"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -1.87%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "1.812.94"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +0.04%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "4.703"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -2.71%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "0.6263"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -1.85%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "2.047.13"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -1.55%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "0.000009189"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -9.12%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "75.27"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -5.36%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "28.586.66"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -2.19%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "5.453"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -7.90%  "

$ws.Range("E20").Value = "  +0.21%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "210.16"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -8.31%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "11.42"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -3.57%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "6.838"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -2.86%  "

$ws.Range("E24").Value = "  +0.18%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "154.41"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -0.06%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "7.837"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -3.26%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "0.1276"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -1.49%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "16.38"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -1.73%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "1.446"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -2.94%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "0.06238"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -4.05%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "1.419"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -2.91%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "3.759"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -2.07%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "3.724"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -2.97%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.706"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -2.26%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.050"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -7.03%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.6352"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -2.95%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "2.507"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -2.11%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "2.716"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -1.56%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.01706"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -2.79%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "6.357"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -2.68%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "1.135.79"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -7.22%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.8644"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -7.27%  "

$ws.Range("E43").Value = "  +0.26%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "100.70"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -0.45%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "1.960.22"
$cell.Style = "Normal"

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "60.47"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -4.82%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.00000000112"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -5.80%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "1.578"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -2.14%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "8.331"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -2.60%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.4500"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -1.60%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.05437"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -1.89%  "
